$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 (copy header formatting from A1, then set the text)
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "link"

# Update row 2 values
$ws.Range("B2").Value = "크리스마스트리 미니트리 풀세트 전구포함 초록트리 트윙클 크리스탈 45cm"
$ws.Range("C2").Value = "https://naver.me/GhS2Sjwk"
$ws.Range("D2").Value = "https://shop-phinf.pstatic.net/20221031_127/1667199162521wbDVQ_JPEG/68335061228336326_657608575.jpg"
$ws.Range("E2").Value = "미니 크리스마스트리 세트"
$ws.Range("F2").Value = "작고 아담한 45cm 트리로 공간을 환하게 꾸며보세요. 전구 포함으로 손쉬운 크리스마스 준비를 도와드립니다."
# G2: empty string, but still present as a real (typed) cell rather than
# being cleared away entirely. Entering a lone apostrophe gives Excel's
# "text" quote-prefix with an empty string; then re-apply the plain
# (non-header) number format from a sibling data cell so no stray
# quote-prefix style lingers on the cell.
$ws.Cells.Item(2, 7).Value = "'"
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)  # xlPasteFormats
